$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 21 new daily rows (2024-08-28 .. 2024-09-25) to the price history table ---

# Column A holds the date as literal text (matches existing rows, which are inlineStr,
# not real Excel dates). Assigning Value2 directly with a "YYYY-MM-DD" string makes Excel
# auto-coerce it into a date serial number, so instead we build each cell as a text formula
# ("="..."") and then flatten the whole A653:A673 block to static values via Copy/PasteSpecial
# (xlPasteValues) -- this avoids picking up any new number-format/style on the cells.

$ws.Cells.Item(653, 1).Formula = '="2024-08-28"'
$ws.Cells.Item(654, 1).Formula = '="2024-08-29"'
$ws.Cells.Item(655, 1).Formula = '="2024-08-30"'
$ws.Cells.Item(656, 1).Formula = '="2024-09-02"'
$ws.Cells.Item(657, 1).Formula = '="2024-09-03"'
$ws.Cells.Item(658, 1).Formula = '="2024-09-04"'
$ws.Cells.Item(659, 1).Formula = '="2024-09-05"'
$ws.Cells.Item(660, 1).Formula = '="2024-09-06"'
$ws.Cells.Item(661, 1).Formula = '="2024-09-09"'
$ws.Cells.Item(662, 1).Formula = '="2024-09-10"'
$ws.Cells.Item(663, 1).Formula = '="2024-09-11"'
$ws.Cells.Item(664, 1).Formula = '="2024-09-12"'
$ws.Cells.Item(665, 1).Formula = '="2024-09-13"'
$ws.Cells.Item(666, 1).Formula = '="2024-09-16"'
$ws.Cells.Item(667, 1).Formula = '="2024-09-17"'
$ws.Cells.Item(668, 1).Formula = '="2024-09-18"'
$ws.Cells.Item(669, 1).Formula = '="2024-09-19"'
$ws.Cells.Item(670, 1).Formula = '="2024-09-20"'
$ws.Cells.Item(671, 1).Formula = '="2024-09-23"'
$ws.Cells.Item(672, 1).Formula = '="2024-09-24"'
$ws.Cells.Item(673, 1).Formula = '="2024-09-25"'

$dateRange = $ws.Range("A653:A673")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Numeric columns: BATA INDIA(C), WHIRLPOOL(D), METROBRANDS(E), CROMPTON(F), VOLTAS(G),
# Basket Value(H), daily return(I), NAV(J)
$ws.Cells.Item(653, 3).Value2 = 1746.650024414062
$ws.Cells.Item(653, 4).Value2 = 1322.349975585938
$ws.Cells.Item(653, 5).Value2 = 1771.75
$ws.Cells.Item(653, 6).Value2 = 2073.949951171875
$ws.Cells.Item(653, 7).Value2 = 462.7999877929688
$ws.Cells.Item(653, 8).Value2 = 8303.099914550781
$ws.Cells.Item(653, 9).Value2 = 0
$ws.Cells.Item(653, 10).Value2 = 236.1929339598256

$ws.Cells.Item(654, 3).Value2 = 1695.900024414062
$ws.Cells.Item(654, 4).Value2 = 1307.849975585938
$ws.Cells.Item(654, 5).Value2 = 1788.400024414062
$ws.Cells.Item(654, 6).Value2 = 2194.800048828125
$ws.Cells.Item(654, 7).Value2 = 465
$ws.Cells.Item(654, 8).Value2 = 8381.950073242188
$ws.Cells.Item(654, 9).Value2 = 0.009496472342001468
$ws.Cells.Item(654, 10).Value2 = 238.4359336245513

$ws.Cells.Item(655, 3).Value2 = 1703.900024414062
$ws.Cells.Item(655, 4).Value2 = 1306.050048828125
$ws.Cells.Item(655, 5).Value2 = 1743.650024414062
$ws.Cells.Item(655, 6).Value2 = 2183.5
$ws.Cells.Item(655, 7).Value2 = 477.0499877929688
$ws.Cells.Item(655, 8).Value2 = 8368.250061035156
$ws.Cells.Item(655, 9).Value2 = -0.001634465975974491
$ws.Cells.Item(655, 10).Value2 = 238.0462182035922

$ws.Cells.Item(656, 3).Value2 = 1660.599975585938
$ws.Cells.Item(656, 4).Value2 = 1280.449951171875
$ws.Cells.Item(656, 5).Value2 = 1769.650024414062
$ws.Cells.Item(656, 6).Value2 = 2220.5
$ws.Cells.Item(656, 7).Value2 = 464.8500061035156
$ws.Cells.Item(656, 8).Value2 = 8325.749969482422
$ws.Cells.Item(656, 9).Value2 = -0.005078731065963999
$ws.Cells.Item(656, 10).Value2 = 236.8372454800664

$ws.Cells.Item(657, 3).Value2 = 1683.75
$ws.Cells.Item(657, 4).Value2 = 1248.449951171875
$ws.Cells.Item(657, 5).Value2 = 1810.949951171875
$ws.Cells.Item(657, 6).Value2 = 2205.5
$ws.Cells.Item(657, 7).Value2 = 469.6000061035156
$ws.Cells.Item(657, 8).Value2 = 8357.449920654297
$ws.Cells.Item(657, 9).Value2 = 0.003807458942205738
$ws.Cells.Item(657, 10).Value2 = 237.7389935682168

$ws.Cells.Item(658, 3).Value2 = 1699.5
$ws.Cells.Item(658, 4).Value2 = 1288.900024414062
$ws.Cells.Item(658, 5).Value2 = 1780.25
$ws.Cells.Item(658, 6).Value2 = 2217.949951171875
$ws.Cells.Item(658, 7).Value2 = 466.5499877929688
$ws.Cells.Item(658, 8).Value2 = 8386.249938964844
$ws.Cells.Item(658, 9).Value2 = 0.003446029420932761
$ws.Cells.Item(658, 10).Value2 = 238.5582491345559

$ws.Cells.Item(659, 3).Value2 = 1694.699951171875
$ws.Cells.Item(659, 4).Value2 = 1272.300048828125
$ws.Cells.Item(659, 5).Value2 = 1783.150024414062
$ws.Cells.Item(659, 6).Value2 = 2236.14990234375
$ws.Cells.Item(659, 7).Value2 = 467.75
$ws.Cells.Item(659, 8).Value2 = 8389.549926757812
$ws.Cells.Item(659, 9).Value2 = 0.0003934998142180441
$ws.Cells.Item(659, 10).Value2 = 238.6521217612705

$ws.Cells.Item(660, 3).Value2 = 1724.449951171875
$ws.Cells.Item(660, 4).Value2 = 1234.300048828125
$ws.Cells.Item(660, 5).Value2 = 1778.650024414062
$ws.Cells.Item(660, 6).Value2 = 2224.39990234375
$ws.Cells.Item(660, 7).Value2 = 463.3999938964844
$ws.Cells.Item(660, 8).Value2 = 8351.999908447266
$ws.Cells.Item(660, 9).Value2 = -0.004475808432915339
$ws.Cells.Item(660, 10).Value2 = 237.5839605821583

$ws.Cells.Item(661, 3).Value2 = 1789.300048828125
$ws.Cells.Item(661, 4).Value2 = 1234.300048828125
$ws.Cells.Item(661, 5).Value2 = 1817.949951171875
$ws.Cells.Item(661, 6).Value2 = 2138.60009765625
$ws.Cells.Item(661, 7).Value2 = 460.7999877929688
$ws.Cells.Item(661, 8).Value2 = 8362.550109863281
$ws.Cells.Item(661, 9).Value2 = 0.001263194627833399
$ws.Cells.Item(661, 10).Value2 = 237.8840753648251

$ws.Cells.Item(662, 3).Value2 = 1799.949951171875
$ws.Cells.Item(662, 4).Value2 = 1249.650024414062
$ws.Cells.Item(662, 5).Value2 = 1830.099975585938
$ws.Cells.Item(662, 6).Value2 = 2140.14990234375
$ws.Cells.Item(662, 7).Value2 = 466.9500122070312
$ws.Cells.Item(662, 8).Value2 = 8420.699890136719
$ws.Cells.Item(662, 9).Value2 = 0.006953594239734628
$ws.Cells.Item(662, 10).Value2 = 239.5382247010065

$ws.Cells.Item(663, 3).Value2 = 1788.25
$ws.Cells.Item(663, 4).Value2 = 1255
$ws.Cells.Item(663, 5).Value2 = 1828.699951171875
$ws.Cells.Item(663, 6).Value2 = 2119.85009765625
$ws.Cells.Item(663, 7).Value2 = 459.2999877929688
$ws.Cells.Item(663, 8).Value2 = 8369.700012207031
$ws.Cells.Item(663, 9).Value2 = -0.006056489198650145
$ws.Cells.Item(663, 10).Value2 = 238.087464030441

$ws.Cells.Item(664, 3).Value2 = 1816.650024414062
$ws.Cells.Item(664, 4).Value2 = 1257.449951171875
$ws.Cells.Item(664, 5).Value2 = 1852.949951171875
$ws.Cells.Item(664, 6).Value2 = 2111
$ws.Cells.Item(664, 7).Value2 = 459
$ws.Cells.Item(664, 8).Value2 = 8415.049926757812
$ws.Cells.Item(664, 9).Value2 = 0.005418344084571652
$ws.Cells.Item(664, 10).Value2 = 239.377503832781

$ws.Cells.Item(665, 3).Value2 = 1888
$ws.Cells.Item(665, 4).Value2 = 1264.349975585938
$ws.Cells.Item(665, 5).Value2 = 1921.550048828125
$ws.Cells.Item(665, 6).Value2 = 2101.35009765625
$ws.Cells.Item(665, 7).Value2 = 451.75
$ws.Cells.Item(665, 8).Value2 = 8530.500122070312
$ws.Cells.Item(665, 9).Value2 = 0.01371949023681921
$ws.Cells.Item(665, 10).Value2 = 242.661641159529

$ws.Cells.Item(666, 3).Value2 = 1934.900024414062
$ws.Cells.Item(666, 4).Value2 = 1259.75
$ws.Cells.Item(666, 5).Value2 = 1916
$ws.Cells.Item(666, 6).Value2 = 2047.199951171875
$ws.Cells.Item(666, 7).Value2 = 447.8999938964844
$ws.Cells.Item(666, 8).Value2 = 8501.54995727539
$ws.Cells.Item(666, 9).Value2 = -0.003393724210849177
$ws.Cells.Item(666, 10).Value2 = 241.8381144728815

$ws.Cells.Item(667, 3).Value2 = 1891.199951171875
$ws.Cells.Item(667, 4).Value2 = 1257.550048828125
$ws.Cells.Item(667, 5).Value2 = 1904.050048828125
$ws.Cells.Item(667, 6).Value2 = 2034.699951171875
$ws.Cells.Item(667, 7).Value2 = 452.75
$ws.Cells.Item(667, 8).Value2 = 8445.75
$ws.Cells.Item(667, 9).Value2 = -0.006563504014657771
$ws.Cells.Item(667, 10).Value2 = 240.2508090376415

$ws.Cells.Item(668, 3).Value2 = 1879.449951171875
$ws.Cells.Item(668, 4).Value2 = 1232.050048828125
$ws.Cells.Item(668, 5).Value2 = 1889.400024414062
$ws.Cells.Item(668, 6).Value2 = 2042.550048828125
$ws.Cells.Item(668, 7).Value2 = 448.8500061035156
$ws.Cells.Item(668, 8).Value2 = 8390.000091552734
$ws.Cells.Item(668, 9).Value2 = -0.006600942302017657
$ws.Cells.Item(668, 10).Value2 = 238.664927309171

$ws.Cells.Item(669, 3).Value2 = 1866.650024414062
$ws.Cells.Item(669, 4).Value2 = 1234.5
$ws.Cells.Item(669, 5).Value2 = 1929.199951171875
$ws.Cells.Item(669, 6).Value2 = 2025.699951171875
$ws.Cells.Item(669, 7).Value2 = 450
$ws.Cells.Item(669, 8).Value2 = 8406.049926757812
$ws.Cells.Item(669, 9).Value2 = 0.00191297199403341
$ws.Cells.Item(669, 10).Value2 = 239.1214866310714

$ws.Cells.Item(670, 3).Value2 = 1930.099975585938
$ws.Cells.Item(670, 4).Value2 = 1246.550048828125
$ws.Cells.Item(670, 5).Value2 = 1928.400024414062
$ws.Cells.Item(670, 6).Value2 = 2062.300048828125
$ws.Cells.Item(670, 7).Value2 = 446.7999877929688
$ws.Cells.Item(670, 8).Value2 = 8507.750061035156
$ws.Cells.Item(670, 9).Value2 = 0.01209844518691423
$ws.Cells.Item(670, 10).Value2 = 242.0144848300909

$ws.Cells.Item(671, 3).Value2 = 2012.849975585938
$ws.Cells.Item(671, 4).Value2 = 1229.900024414062
$ws.Cells.Item(671, 5).Value2 = 1882.449951171875
$ws.Cells.Item(671, 6).Value2 = 2039
$ws.Cells.Item(671, 7).Value2 = 445.75
$ws.Cells.Item(671, 8).Value2 = 8501.449951171875
$ws.Cells.Item(671, 9).Value2 = -0.000740514215636784
$ws.Cells.Item(671, 10).Value2 = 241.8352696636842

$ws.Cells.Item(672, 3).Value2 = 2026
$ws.Cells.Item(672, 4).Value2 = 1238.800048828125
$ws.Cells.Item(672, 5).Value2 = 1914.400024414062
$ws.Cells.Item(672, 6).Value2 = 2064.550048828125
$ws.Cells.Item(672, 7).Value2 = 439.25
$ws.Cells.Item(672, 8).Value2 = 8561.500122070312
$ws.Cells.Item(672, 9).Value2 = 0.007063521075032611
$ws.Cells.Item(672, 10).Value2 = 243.5434781876398

$ws.Cells.Item(673, 3).Value2 = 2051.39990234375
$ws.Cells.Item(673, 4).Value2 = 1252.949951171875
$ws.Cells.Item(673, 5).Value2 = 1910.849975585938
$ws.Cells.Item(673, 6).Value2 = 2112.050048828125
$ws.Cells.Item(673, 7).Value2 = 441.5499877929688
$ws.Cells.Item(673, 8).Value2 = 8651.899841308594
$ws.Cells.Item(673, 9).Value2 = 0.01055886444540762
$ws.Cells.Item(673, 10).Value2 = 246.1150207603862

